$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 240.76923
$ws.Range("I33").Value = 235.91667
$ws.Range("J33").Value = 299
$ws.Range("K33").Value = 235.91667
$ws.Range("L33").Value = 299
$ws.Range("M33").Value = -6.916670000000011
$ws.Range("N33").Value = -757
$ws.Range("H88").Value = 6006.3
$ws.Range("I88").Value = 2949.5
$ws.Range("J88").Value = 6770.5
$ws.Range("K88").Value = 2949.5
$ws.Range("L88").Value = 6770.5
$ws.Range("M88").Value = -2543.5
$ws.Range("N88").Value = -7582.5
$ws.Range("H91").Value = 6006.3
$ws.Range("I91").Value = 2949.5
$ws.Range("J91").Value = 6770.5
$ws.Range("K91").Value = 2949.5
$ws.Range("L91").Value = 6770.5
$ws.Range("M91").Value = -1545.5
$ws.Range("N91").Value = -9578.5
$ws.Range("H138").Value = 7372.7954
$ws.Range("I138").Value = 11046.521
$ws.Range("K138").Value = 33139.563
$ws.Range("M138").Value = -27999.563

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1520.2059
$ws.Range("I32").Value = 1365.2188
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 1365.2188
$ws.Range("L32").Value = 4000
$ws.Range("M32").Value = -1078.2188
$ws.Range("N32").Value = -4574
$ws.Range("H97").Value = 673.1905
$ws.Range("I97").Value = 673.1905
$ws.Range("K97").Value = 673.1905
$ws.Range("M97").Value = -177.1905
$ws.Range("H120").Value = 77946.336
$ws.Range("J120").Value = 77946.336
$ws.Range("L120").Value = 77946.336
$ws.Range("N120").Value = -87622.336
$ws.Range("H122").Value = 2461.5
$ws.Range("I122").Value = 1697
$ws.Range("J122").Value = 3365
$ws.Range("K122").Value = 5091
$ws.Range("L122").Value = 10095
$ws.Range("M122").Value = -2641
$ws.Range("N122").Value = -14995
$ws.Range("H132").Value = 2866.1707
$ws.Range("I132").Value = 1891.2258
$ws.Range("K132").Value = 5673.6774
$ws.Range("M132").Value = -3143.6774
$ws.Range("H139").Value = 90053.5
$ws.Range("J139").Value = 90053.5
$ws.Range("L139").Value = 90053.5
$ws.Range("N139").Value = -100333.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1812.5
$ws.Range("I86").Value = 1689.9333
$ws.Range("J86").Value = 2180.2
$ws.Range("K86").Value = 1689.9333
$ws.Range("L86").Value = 2180.2
$ws.Range("M86").Value = -566.9332999999999
$ws.Range("N86").Value = -4426.2
$ws.Range("H89").Value = 1812.5
$ws.Range("I89").Value = 1689.9333
$ws.Range("J89").Value = 2180.2
$ws.Range("K89").Value = 8449.666499999999
$ws.Range("L89").Value = 10901
$ws.Range("M89").Value = -2833.666499999999
$ws.Range("N89").Value = -22133

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9170.48
$ws.Range("I31").Value = 1455.9656
$ws.Range("J31").Value = 18897.479
$ws.Range("K31").Value = 1455.9656
$ws.Range("L31").Value = 18897.479
$ws.Range("M31").Value = -1160.9656
$ws.Range("N31").Value = -19487.479
$ws.Range("H34").Value = 9170.48
$ws.Range("I34").Value = 1455.9656
$ws.Range("J34").Value = 18897.479
$ws.Range("K34").Value = 1455.9656
$ws.Range("L34").Value = 18897.479
$ws.Range("M34").Value = -1253.9656
$ws.Range("N34").Value = -19301.479
$ws.Range("H62").Value = 8112.1
$ws.Range("J62").Value = 10444.857
$ws.Range("L62").Value = 10444.857
$ws.Range("N62").Value = -11692.857
$ws.Range("H65").Value = 8112.1
$ws.Range("J65").Value = 10444.857
$ws.Range("L65").Value = 52224.285
$ws.Range("N65").Value = -58464.285
$ws.Range("H86").Value = 9217.25
$ws.Range("J86").Value = 9949.5
$ws.Range("L86").Value = 9949.5
$ws.Range("N86").Value = -12195.5
$ws.Range("H88").Value = 45486.75
$ws.Range("J88").Value = 49725
$ws.Range("L88").Value = 49725
$ws.Range("N88").Value = -50537
$ws.Range("H89").Value = 9217.25
$ws.Range("J89").Value = 9949.5
$ws.Range("L89").Value = 49747.5
$ws.Range("N89").Value = -60979.5
$ws.Range("H91").Value = 45486.75
$ws.Range("J91").Value = 49725
$ws.Range("L91").Value = 49725
$ws.Range("N91").Value = -52533
$ws.Range("H134").Value = 2602.1738
$ws.Range("I134").Value = 2096.2354
$ws.Range("K134").Value = 6288.706200000001
$ws.Range("M134").Value = -3753.706200000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10039.615
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 7500
$ws.Range("M68").Value = -6689
$ws.Range("H71").Value = 10039.615
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 22500
$ws.Range("M71").Value = -18444
$ws.Range("H132").Value = 1330.3334
$ws.Range("I132").Value = 936.75
$ws.Range("J132").Value = 2117.5
$ws.Range("K132").Value = 8430.75
$ws.Range("L132").Value = 19057.5
$ws.Range("M132").Value = -5900.75
$ws.Range("N132").Value = -24117.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2634.625
$ws.Range("I102").Value = 2518.2144
$ws.Range("K102").Value = 2518.2144
$ws.Range("M102").Value = -896.2143999999998
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 47607.69
$ws.Range("I122").Value = 75492.78999999999
$ws.Range("J122").Value = 15075.083
$ws.Range("K122").Value = 226478.37
$ws.Range("L122").Value = 45225.249
$ws.Range("M122").Value = -224028.37
$ws.Range("N122").Value = -50125.249
$ws.Range("H132").Value = 11633416
$ws.Range("I132").Value = 16672122
$ws.Range("J132").Value = 5632.769
$ws.Range("K132").Value = 50016366
$ws.Range("L132").Value = 16898.307
$ws.Range("M132").Value = -50013836
$ws.Range("N132").Value = -21958.307

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1127.4166
$ws.Range("I22").Value = 407.33334
$ws.Range("J22").Value = 1847.5
$ws.Range("K22").Value = 407.33334
$ws.Range("L22").Value = 1847.5
$ws.Range("M22").Value = -112.33334
$ws.Range("N22").Value = -2437.5
$ws.Range("H27").Value = 1127.4166
$ws.Range("I27").Value = 407.33334
$ws.Range("J27").Value = 1847.5
$ws.Range("K27").Value = 407.33334
$ws.Range("L27").Value = 1847.5
$ws.Range("M27").Value = -300.33334
$ws.Range("N27").Value = -2061.5
$ws.Range("H61").Value = 1344.0769
$ws.Range("I61").Value = 1225.4546
$ws.Range("K61").Value = 1225.4546
$ws.Range("M61").Value = -1023.4546
$ws.Range("H113").Value = 1344.0769
$ws.Range("I113").Value = 1225.4546
$ws.Range("K113").Value = 1225.4546
$ws.Range("M113").Value = 944.5454
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2280.7273
$ws.Range("J126").Value = 4994.6665
$ws.Range("L126").Value = 14983.9995
$ws.Range("N126").Value = -19923.9995
